$wb = $excel.ActiveWorkbook

# ----- Sheet: Cases -----
$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X65").Value = 5273
$ws.Range("AB65").Value = 29097
$ws.Range("X66").Value = 5283
$ws.Range("AB66").Value = 29203
$ws.Range("X67").Value = 5284
$ws.Range("AB67").Value = 29328
$ws.Range("AB68").Value = 29367

# ----- Sheet: Fatalities -----
$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("X66").Value = 376
$ws.Range("AB66").Value = 1730
$ws.Range("X67").Value = 377
$ws.Range("AB67").Value = 1740
$ws.Range("AB68").Value = 1743

# ----- Sheet: Hospitalized -----
$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X27").Value = 210
$ws.Range("AB27").Value = 904
$ws.Range("X28").Value = 230
$ws.Range("AB28").Value = 1090
$ws.Range("X29").Value = 248
$ws.Range("AB29").Value = 1207
$ws.Range("X30").Value = 284
$ws.Range("AB30").Value = 1364
$ws.Range("X35").Value = 366
$ws.Range("AB35").Value = 1990
$ws.Range("X36").Value = 368
$ws.Range("AB36").Value = 2174
$ws.Range("X37").Value = 376
$ws.Range("AB37").Value = 2202
$ws.Range("X38").Value = 380
$ws.Range("AB38").Value = 2280
$ws.Range("X39").Value = 370
$ws.Range("AB39").Value = 2343
$ws.Range("X40").Value = 361
$ws.Range("AB40").Value = 2324
$ws.Range("X41").Value = 360
$ws.Range("AB41").Value = 2311
$ws.Range("X42").Value = 371
$ws.Range("AB42").Value = 2296
$ws.Range("X43").Value = 356
$ws.Range("AB43").Value = 2301
$ws.Range("X44").Value = 329
$ws.Range("AB44").Value = 2219
$ws.Range("X45").Value = 320
$ws.Range("AB45").Value = 2129
$ws.Range("X46").Value = 303
$ws.Range("AB46").Value = 2061
$ws.Range("X47").Value = 294
$ws.Range("AB47").Value = 2001
$ws.Range("X48").Value = 280
$ws.Range("AB48").Value = 1926
$ws.Range("X49").Value = 289
$ws.Range("AB49").Value = 1903
$ws.Range("X50").Value = 289
$ws.Range("AB50").Value = 1886
$ws.Range("X51").Value = 272
$ws.Range("AB51").Value = 1846
$ws.Range("X52").Value = 245
$ws.Range("AB52").Value = 1723
$ws.Range("X53").Value = 236
$ws.Range("AB53").Value = 1665
$ws.Range("X54").Value = 220
$ws.Range("AB54").Value = 1566
$ws.Range("X55").Value = 217
$ws.Range("AB55").Value = 1515
$ws.Range("X56").Value = 218
$ws.Range("AB56").Value = 1507
$ws.Range("X57").Value = 201
$ws.Range("AB57").Value = 1483
$ws.Range("X58").Value = 188
$ws.Range("AB58").Value = 1405
$ws.Range("X59").Value = 174
$ws.Range("AB59").Value = 1339
$ws.Range("X60").Value = 174
$ws.Range("AB60").Value = 1286
$ws.Range("X61").Value = 162
$ws.Range("AB61").Value = 1241
$ws.Range("X62").Value = 167
$ws.Range("AB62").Value = 1212
$ws.Range("X63").Value = 175
$ws.Range("AB63").Value = 1191
$ws.Range("X64").Value = 166
$ws.Range("AB64").Value = 1169
$ws.Range("X65").Value = 156
$ws.Range("AB65").Value = 1141
$ws.Range("X66").Value = 146
$ws.Range("AB66").Value = 1072
$ws.Range("X67").Value = 142
$ws.Range("AB67").Value = 998
$ws.Range("AB68").Value = 992

# ----- Sheet: ICU -----
$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X38").Value = 67
$ws.Range("AB38").Value = 383
$ws.Range("X39").Value = 68
$ws.Range("AB39").Value = 389
$ws.Range("X40").Value = 66
$ws.Range("AB40").Value = 395
$ws.Range("X41").Value = 66
$ws.Range("AB41").Value = 395
$ws.Range("X42").Value = 65
$ws.Range("AB42").Value = 385
$ws.Range("X43").Value = 63
$ws.Range("AB43").Value = 385
$ws.Range("X44").Value = 62
$ws.Range("AB44").Value = 377
$ws.Range("X45").Value = 61
$ws.Range("AB45").Value = 378
$ws.Range("X46").Value = 61
$ws.Range("AB46").Value = 371
$ws.Range("X47").Value = 63
$ws.Range("AB47").Value = 364
$ws.Range("X48").Value = 63
$ws.Range("AB48").Value = 364
$ws.Range("X49").Value = 62
$ws.Range("AB49").Value = 354
$ws.Range("X50").Value = 58
$ws.Range("AB50").Value = 349
$ws.Range("X51").Value = 57
$ws.Range("AB51").Value = 334
$ws.Range("X52").Value = 53
$ws.Range("AB52").Value = 313
$ws.Range("X53").Value = 52
$ws.Range("AB53").Value = 292
$ws.Range("X54").Value = 51
$ws.Range("AB54").Value = 286
$ws.Range("X55").Value = 48
$ws.Range("AB55").Value = 273
$ws.Range("X56").Value = 47
$ws.Range("AB56").Value = 265
$ws.Range("X57").Value = 41
$ws.Range("AB57").Value = 251
$ws.Range("X58").Value = 41
$ws.Range("AB58").Value = 239
$ws.Range("X59").Value = 40
$ws.Range("AB59").Value = 205
$ws.Range("X60").Value = 38
$ws.Range("AB60").Value = 197
$ws.Range("X61").Value = 35
$ws.Range("AB61").Value = 195
$ws.Range("X62").Value = 35
$ws.Range("AB62").Value = 190
$ws.Range("X63").Value = 35
$ws.Range("AB63").Value = 180
$ws.Range("X64").Value = 33
$ws.Range("AB64").Value = 177
$ws.Range("X65").Value = 32
$ws.Range("AB65").Value = 163
$ws.Range("X66").Value = 31
$ws.Range("AB66").Value = 168
$ws.Range("X67").Value = 27
$ws.Range("AB67").Value = 157
$ws.Range("AB68").Value = 154
